$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 230 (pushes the existing rows 230-232 down
# to 231-233), then populate it with the new "bead" term.
$ws.Rows.Item(230).Insert()

# Set label ("bead") before the IRI so the shared-string table gets the
# same ordering as the authored workbook (bead -> idx 638, IRI -> idx 639).
$ws.Cells.Item(230, 2).Value = "bead"
$ws.Cells.Item(230, 1).Value = "http://purl.obolibrary.org/obo/OBI_1000207"
$ws.Cells.Item(230, 3).Value = "y"

# Match the column C formatting used by the rest of the "Include in View"
# column (style index 4 in the original workbook) by copying the format
# from the neighboring row.
$ws.Cells.Item(231, 3).Copy()
$ws.Cells.Item(230, 3).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the view state (scroll position / selection) recorded for the
# sheet after the edit.
$ws.Range("C237").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 223
$win.ScrollColumn = 1
